$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = -7.458299999999996
    "C3"  = -11.5533
    "C14" = -13.13919999999999
    "C21" = -12.4191
    "C23" = -12.58350000000001
    "C25" = -14.25479999999999
    "D25" = -8.048
    "C26" = -12.80270000000001
    "D27" = -8.710700000000005
    "C29" = -10.82660000000001
    "D31" = -8.560300000000009
    "D39" = -8.122700000000002
    "D48" = -7.327699999999998
    "D51" = -7.723499999999996
    "D52" = -7.873799999999996
    "C53" = -10.48180000000001
    "D55" = -8.321399999999999
    "D56" = -7.934299999999999
    "C57" = -14.02569999999999
    "D57" = -8.2987
    "C59" = -12.66859999999999
    "C69" = -10.914
    "D73" = -7.798699999999999
    "C79" = -10.77600000000002
    "C83" = -13.74949999999999
    "D89" = -6.082800000000001
    "D90" = -8.120300000000004
    "C91" = -10.3765
    "D92" = -6.337300000000002
    "C93" = -11.4566
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
